$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the header of column C (row 2): "Titúlo de la idea" -> "Título de la idea"
$ws.Range("C2").Value = "Título de la idea"

# Update the active selection to C2 (as reflected in the saved file)
$ws.Range("C2").Select()
